$wb = $excel.ActiveWorkbook

# ALC!row40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2524.4167
$ws.Range("I40").Value = 1438.2
$ws.Range("J40").Value = 3300.2856
$ws.Range("K40").Value = 1438.2
$ws.Range("L40").Value = 3300.2856
$ws.Range("M40").Value = -1263.2
$ws.Range("N40").Value = -3650.2856

# ALC!row43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1552.2222
$ws.Range("I43").Value = 978
$ws.Range("J43").Value = 2270
$ws.Range("K43").Value = 978
$ws.Range("L43").Value = 2270
$ws.Range("M43").Value = -909
$ws.Range("N43").Value = -2408

# ALC!row76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2895.6287
$ws.Range("I76").Value = 2618.5715
$ws.Range("J76").Value = 3542.0952
$ws.Range("K76").Value = 2618.5715
$ws.Range("L76").Value = 3542.0952
$ws.Range("M76").Value = -2303.5715
$ws.Range("N76").Value = -4172.0952

# ALC!row79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 2895.6287
$ws.Range("I79").Value = 2618.5715
$ws.Range("J79").Value = 3542.0952
$ws.Range("K79").Value = 2618.5715
$ws.Range("L79").Value = 3542.0952
$ws.Range("M79").Value = -1526.5715
$ws.Range("N79").Value = -5726.0952

# ALC!row86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1770
$ws.Range("I86").Value = 1762.5
$ws.Range("K86").Value = 1762.5
$ws.Range("M86").Value = -639.5

# ALC!row89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1770
$ws.Range("I89").Value = 1762.5
$ws.Range("K89").Value = 8812.5
$ws.Range("M89").Value = -3196.5

# ARM!row107
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 20000
$ws.Range("J107").Value = 20000
$ws.Range("L107").Value = 20000
$ws.Range("N107").Value = -27680

# BSM!row20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2125.6
$ws.Range("I20").Value = 1840.9
$ws.Range("J20").Value = 2505.2
$ws.Range("K20").Value = 1840.9
$ws.Range("L20").Value = 2505.2
$ws.Range("M20").Value = -1593.9
$ws.Range("N20").Value = -2999.2

# BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1834.3636
$ws.Range("I86").Value = 2350
$ws.Range("J86").Value = 1719.7778
$ws.Range("K86").Value = 2350
$ws.Range("L86").Value = 1719.7778
$ws.Range("M86").Value = -1227
$ws.Range("N86").Value = -3965.7778

# BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1834.3636
$ws.Range("I89").Value = 2350
$ws.Range("J89").Value = 1719.7778
$ws.Range("K89").Value = 11750
$ws.Range("L89").Value = 8598.889000000001
$ws.Range("M89").Value = -6134
$ws.Range("N89").Value = -19830.889

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 693.4
$ws.Range("I94").Value = 593.5
$ws.Range("J94").Value = 760
$ws.Range("K94").Value = 593.5
$ws.Range("L94").Value = 760
$ws.Range("M94").Value = -142.5
$ws.Range("N94").Value = -1662

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1239.3695
$ws.Range("I99").Value = 1039.963
$ws.Range("J99").Value = 1522.7368
$ws.Range("K99").Value = 1039.963
$ws.Range("L99").Value = 1522.7368
$ws.Range("M99").Value = 458.037
$ws.Range("N99").Value = -4518.7368

# BSM!row112
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 29750
$ws.Range("J112").Value = 29750
$ws.Range("L112").Value = 29750
$ws.Range("N112").Value = -32704

# CUL!row51
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3571.4285
$ws.Range("I51").Value = 1000
$ws.Range("K51").Value = 3000
$ws.Range("M51").Value = -2540

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 932666
$ws.Range("I122").Value = 6626.7427
$ws.Range("J122").Value = 3247764
$ws.Range("K122").Value = 59640.6843
$ws.Range("L122").Value = 29229876
$ws.Range("M122").Value = -57190.6843
$ws.Range("N122").Value = -29234776

# CUL!row129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4423.3945
$ws.Range("I129").Value = 1597.8
$ws.Range("J129").Value = 6266.174
$ws.Range("K129").Value = 4793.4
$ws.Range("L129").Value = 18798.522
$ws.Range("M129").Value = 206.6000000000004
$ws.Range("N129").Value = -28798.522

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 780.13794
$ws.Range("I131").Value = 421.25
$ws.Range("J131").Value = 916.8570999999999
$ws.Range("K131").Value = 1263.75
$ws.Range("L131").Value = 2750.5713
$ws.Range("M131").Value = 3776.25
$ws.Range("N131").Value = -12830.5713

# CUL!row133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4086.3635
$ws.Range("I133").Value = 2992.8572
$ws.Range("J133").Value = 6000
$ws.Range("K133").Value = 8978.571599999999
$ws.Range("L133").Value = 18000
$ws.Range("M133").Value = -3918.571599999999
$ws.Range("N133").Value = -28120

# CUL!row140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4487.4736
$ws.Range("I140").Value = 3191.8823
$ws.Range("J140").Value = 15500
$ws.Range("K140").Value = 9575.6469
$ws.Range("L140").Value = 46500
$ws.Range("M140").Value = -4395.6469
$ws.Range("N140").Value = -56860

# GSM!row2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 42.38889
$ws.Range("I2").Value = 62.7
$ws.Range("K2").Value = 62.7
$ws.Range("M2").Value = 50.3

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4640.6113
$ws.Range("I132").Value = 4610.8184
$ws.Range("J132").Value = 4687.4287
$ws.Range("K132").Value = 13832.4552
$ws.Range("L132").Value = 14062.2861
$ws.Range("M132").Value = -11302.4552
$ws.Range("N132").Value = -19122.2861

# GSM!row133
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 39944
$ws.Range("J133").Value = 39944
$ws.Range("L133").Value = 39944
$ws.Range("N133").Value = -50064

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1015.5455
$ws.Range("I22").Value = 638.2
$ws.Range("J22").Value = 1330
$ws.Range("K22").Value = 638.2
$ws.Range("L22").Value = 1330
$ws.Range("M22").Value = -343.2
$ws.Range("N22").Value = -1920

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1015.5455
$ws.Range("I27").Value = 638.2
$ws.Range("J27").Value = 1330
$ws.Range("K27").Value = 638.2
$ws.Range("L27").Value = 1330
$ws.Range("M27").Value = -531.2
$ws.Range("N27").Value = -1544

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1057.7778
$ws.Range("I46").Value = 596.6667
$ws.Range("J46").Value = 1288.3334
$ws.Range("K46").Value = 596.6667
$ws.Range("L46").Value = 1288.3334
$ws.Range("M46").Value = -408.6667
$ws.Range("N46").Value = -1664.3334

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2811.8823
$ws.Range("I68").Value = 2781.375
$ws.Range("K68").Value = 2781.375
$ws.Range("M68").Value = -2032.375

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2811.8823
$ws.Range("I71").Value = 2781.375
$ws.Range("K71").Value = 13906.875
$ws.Range("M71").Value = -10162.875

# LTW!row100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 20001698
$ws.Range("I100").Value = 1705.1875
$ws.Range("J100").Value = 55557240
$ws.Range("K100").Value = 1705.1875
$ws.Range("L100").Value = 55557240
$ws.Range("M100").Value = -1164.1875
$ws.Range("N100").Value = -55558322

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3055.9355
$ws.Range("I132").Value = 2202.875
$ws.Range("J132").Value = 3965.8667
$ws.Range("K132").Value = 6608.625
$ws.Range("L132").Value = 11897.6001
$ws.Range("M132").Value = -4078.625
$ws.Range("N132").Value = -16957.6001
